$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.407.14'
$ws.Range("E2").Value = '  +1.61%  '

$ws.Range("D3").Value = '2.237.31'
$ws.Range("E3").Value = '  +0.64%  '

$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").Value = '317.05'
$ws.Range("E5").Value = '  +0.85%  '

$ws.Range("D6").Value = '99.72'
$ws.Range("E6").Value = '  +1.16%  '

$ws.Range("D7").Value = '0.582'
$ws.Range("E7").Value = '  +1.61%  '

$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").Value = '0.563'
$ws.Range("E9").Value = '  +0.74%  '

$ws.Range("D10").Value = '37.11'
$ws.Range("E10").Value = '  +0.16%  '

$ws.Range("D11").Value = '0.0832'
$ws.Range("E11").Value = '  -0.11%  '

$ws.Range("D12").Value = '7.71'
$ws.Range("E12").Value = '  +1.13%  '

$ws.Range("E13").Value = '  +2.06%  '

$ws.Range("D14").Value = '0.865'
$ws.Range("E14").Value = '  -0.88%  '

$ws.Range("D15").Value = '14.34'
$ws.Range("E15").Value = '  +2.77%  '

$ws.Range("D16").Value = '2.249.35'
$ws.Range("E16").Value = '  +1.57%  '

$ws.Range("D17").Value = '43.341.13'
$ws.Range("E17").Value = '  +1.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.67%  '

$ws.Range("E19").Value = '  +0.34%  '

$ws.Range("D20").Value = '0.0₃0974'
$ws.Range("E20").Value = '  +2.66%  '

$ws.Range("D21").Value = '65.37'
$ws.Range("E21").Value = '  +0.95%  '

$ws.Range("E22").Value = '  -2.71%  '

$ws.Range("D23").Value = '236.45'
$ws.Range("E23").Value = '  +0.73%  '

$ws.Range("E24").Value = '  +2.58%  '

$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("D26").Value = '4.05'
$ws.Range("E26").Value = '  +2.97%  '

$ws.Range("D27").Value = '10.06'
$ws.Range("E27").Value = '  -1.14%  '

$ws.Range("E28").Value = '  +2.03%  '

$ws.Range("D29").Value = '6.39'
$ws.Range("E29").Value = '  -2.20%  '

$ws.Range("D30").Value = '36.61'
$ws.Range("E30").Value = '  +11.34%  '

$ws.Range("D31").Value = '20.31'
$ws.Range("E31").Value = '  -0.57%  '

$ws.Range("D32").Value = '0.0871'
$ws.Range("E32").Value = '  -1.86%  '

$ws.Range("D33").Value = '157.74'
$ws.Range("E33").Value = '  -0.64%  '

$ws.Range("D34").Value = '2.71'
$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("E35").Value = '  +3.27%  '

$ws.Range("E36").Value = '  -1.42%  '

$ws.Range("E37").Value = '  +2.37%  '

$ws.Range("D38").Value = '4.38'
$ws.Range("E38").Value = '  -1.69%  '

$ws.Range("E39").Value = '  +0.34%  '

$ws.Range("D40").Value = '3.71'
$ws.Range("E40").Value = '  +4.42%  '

$ws.Range("D41").Value = '0.0322'
$ws.Range("E41").Value = '  -0.28%  '

$ws.Range("D42").Value = '14.45'
$ws.Range("E42").Value = '  +20.25%  '

$ws.Range("E43").Value = '  -0.20%  '

$ws.Range("D44").Value = '1.833.88'
$ws.Range("E44").Value = '  +1.19%  '

$ws.Range("E45").Value = '  -1.81%  '

$ws.Range("D46").Value = '84.34'
$ws.Range("E46").Value = '  -5.37%  '

$ws.Range("E47").Value = '  -1.58%  '

$ws.Range("E48").Value = '  +2.78%  '

$ws.Range("D49").Value = '74.08'
$ws.Range("E49").Value = '  -4.32%  '

$ws.Range("D50").Value = '103.38'
$ws.Range("E50").Value = '  +1.56%  '

$ws.Range("D51").Value = '58.28'
$ws.Range("E51").Value = '  -3.52%  '
